$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New week 19 data (row 20)
$ws.Range("A20").Value = 19

$ws.Range("B20").Value = 2.0889930555555556
$ws.Range("B20").NumberFormat = $ws.Range("B19").NumberFormat

$ws.Range("C20").Formula = "=SUM(B2:B20)+1.2708333333"
$ws.Range("C20").NumberFormat = $ws.Range("C19").NumberFormat

$ws.Range("D20").Value = "Élite (Subtitled, Spanish, New):37; the feynman lectures on physics (Text with visuals, English, New):44; Historia de un crimen: Búsqueda (Subtitled, Spanish, New):39; La casa de papel (Subtitled, Spanish, New):36;"

$ws.Range("D20").Select()
